# Auto-generated Excel COM-interop script
# Applies the scheduled market-data refresh described in the commit:
#   "chore: update Sheets via scheduled runner"
#
# For each affected leve row, refreshes the derived market-price columns
# (currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ) with freshly
# pulled values. Where a profit figure is no longer computable for a row
# (e.g. WVR row 139's HQ profit), the cell is cleared entirely rather than
# left with a stale number, matching upstream's sparse-cell convention.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 709.8333
$ws.Range("I38").Value = 709.8333
$ws.Range("K38").Value = 2129.4999
$ws.Range("M38").Value = -1757.4999
$ws.Range("H43").Value = 1130.7858
$ws.Range("J43").Value = 1130.7858
$ws.Range("L43").Value = 1130.7858
$ws.Range("N43").Value = -1268.7858
$ws.Range("H86").Value = 138735.78
$ws.Range("I86").Value = 177803.14
$ws.Range("K86").Value = 177803.14
$ws.Range("M86").Value = -176680.14
$ws.Range("H89").Value = 138735.78
$ws.Range("I89").Value = 177803.14
$ws.Range("K89").Value = 889015.7000000001
$ws.Range("M89").Value = -883399.7000000001
$ws.Range("H121").Value = 829.2
$ws.Range("J121").Value = 998.6667
$ws.Range("L121").Value = 2996.0001
$ws.Range("N121").Value = -6490.0001
$ws.Range("H138").Value = 1582.3064
$ws.Range("I138").Value = 1348.2727
$ws.Range("J138").Value = 3421.1428
$ws.Range("K138").Value = 4044.8181
$ws.Range("L138").Value = 10263.4284
$ws.Range("M138").Value = 1095.1819
$ws.Range("N138").Value = -20543.4284

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1450
$ws.Range("J4").Value = 1450
$ws.Range("L4").Value = 1450
$ws.Range("N4").Value = -1682
$ws.Range("H32").Value = 3828.0757
$ws.Range("I32").Value = 3084.638
$ws.Range("J32").Value = 9218
$ws.Range("K32").Value = 3084.638
$ws.Range("L32").Value = 9218
$ws.Range("M32").Value = -2797.638
$ws.Range("N32").Value = -9792
$ws.Range("H61").Value = 2227.8096
$ws.Range("I61").Value = 1205.375
$ws.Range("J61").Value = 5499.6
$ws.Range("K61").Value = 1205.375
$ws.Range("L61").Value = 5499.6
$ws.Range("M61").Value = -993.375
$ws.Range("N61").Value = -5923.6
$ws.Range("H74").Value = 1338.8292
$ws.Range("I74").Value = 976.8276
$ws.Range("J74").Value = 2213.6667
$ws.Range("K74").Value = 976.8276
$ws.Range("L74").Value = 2213.6667
$ws.Range("M74").Value = -102.8276
$ws.Range("N74").Value = -3961.6667
$ws.Range("H77").Value = 1338.8292
$ws.Range("I77").Value = 976.8276
$ws.Range("J77").Value = 2213.6667
$ws.Range("K77").Value = 4884.138
$ws.Range("L77").Value = 11068.3335
$ws.Range("M77").Value = -516.1379999999999
$ws.Range("N77").Value = -19804.3335
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("M126").Value = -3530
$ws.Range("H132").Value = 1463.5
$ws.Range("I132").Value = 1100.1
$ws.Range("J132").Value = 2008.6
$ws.Range("K132").Value = 3300.3
$ws.Range("L132").Value = 6025.799999999999
$ws.Range("M132").Value = -770.2999999999997
$ws.Range("N132").Value = -11085.8
$ws.Range("H136").Value = 2227.8096
$ws.Range("I136").Value = 1205.375
$ws.Range("J136").Value = 5499.6
$ws.Range("K136").Value = 3616.125
$ws.Range("L136").Value = 16498.8
$ws.Range("M136").Value = -1066.125
$ws.Range("N136").Value = -21598.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 397.5
$ws.Range("J22").Value = 397.5
$ws.Range("L22").Value = 397.5
$ws.Range("N22").Value = -743.5
$ws.Range("H92").Value = 19499.25
$ws.Range("J92").Value = 19499.25
$ws.Range("L92").Value = 19499.25
$ws.Range("N92").Value = -24491.25
$ws.Range("H134").Value = 5202.4185
$ws.Range("I134").Value = 5002.3687
$ws.Range("K134").Value = 15007.1061
$ws.Range("M134").Value = -12472.1061

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1992.5217
$ws.Range("I31").Value = 1571.8
$ws.Range("J31").Value = 2316.1538
$ws.Range("K31").Value = 1571.8
$ws.Range("L31").Value = 2316.1538
$ws.Range("M31").Value = -1276.8
$ws.Range("N31").Value = -2906.1538
$ws.Range("H34").Value = 1992.5217
$ws.Range("I34").Value = 1571.8
$ws.Range("J34").Value = 2316.1538
$ws.Range("K34").Value = 1571.8
$ws.Range("L34").Value = 2316.1538
$ws.Range("M34").Value = -1369.8
$ws.Range("N34").Value = -2720.1538
$ws.Range("H58").Value = 3109597.8
$ws.Range("I58").Value = 10872540
$ws.Range("J58").Value = 4420.9
$ws.Range("K58").Value = 10872540
$ws.Range("L58").Value = 4420.9
$ws.Range("M58").Value = -10872337
$ws.Range("N58").Value = -4826.9
$ws.Range("H86").Value = 166668670
$ws.Range("I86").Value = 200001540
$ws.Range("K86").Value = 200001540
$ws.Range("M86").Value = -200000417
$ws.Range("H88").Value = 45333.332
$ws.Range("J88").Value = 45333.332
$ws.Range("L88").Value = 45333.332
$ws.Range("N88").Value = -46145.332
$ws.Range("H89").Value = 166668670
$ws.Range("I89").Value = 200001540
$ws.Range("K89").Value = 1000007700
$ws.Range("M89").Value = -1000002084
$ws.Range("H91").Value = 45333.332
$ws.Range("J91").Value = 45333.332
$ws.Range("L91").Value = 45333.332
$ws.Range("N91").Value = -48141.332
$ws.Range("H132").Value = 2213.52
$ws.Range("I132").Value = 1178.8422
$ws.Range("J132").Value = 5490
$ws.Range("K132").Value = 3536.5266
$ws.Range("L132").Value = 16470
$ws.Range("M132").Value = -1006.5266
$ws.Range("N132").Value = -21530
$ws.Range("H134").Value = 1562.1837
$ws.Range("I134").Value = 1475.1351
$ws.Range("J134").Value = 1830.5834
$ws.Range("K134").Value = 4425.4053
$ws.Range("L134").Value = 5491.7502
$ws.Range("M134").Value = -1890.4053
$ws.Range("N134").Value = -10561.7502
$ws.Range("H136").Value = 3109597.8
$ws.Range("I136").Value = 10872540
$ws.Range("J136").Value = 4420.9
$ws.Range("K136").Value = 32617620
$ws.Range("L136").Value = 13262.7
$ws.Range("M136").Value = -32615070
$ws.Range("N136").Value = -18362.7
$ws.Range("H141").Value = 84970
$ws.Range("J141").Value = 84970
$ws.Range("L141").Value = 84970
$ws.Range("N141").Value = -95330

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 651.53845
$ws.Range("I5").Value = 541.44446
$ws.Range("J5").Value = 899.25
$ws.Range("K5").Value = 1624.33338
$ws.Range("L5").Value = 2697.75
$ws.Range("M5").Value = -1512.33338
$ws.Range("N5").Value = -2921.75
$ws.Range("H109").Value = 3846.348
$ws.Range("I109").Value = 1130.5
$ws.Range("J109").Value = 5294.8
$ws.Range("K109").Value = 3391.5
$ws.Range("L109").Value = 15884.4
$ws.Range("M109").Value = -2351.5
$ws.Range("N109").Value = -17964.4
$ws.Range("H122").Value = 981.7222
$ws.Range("I122").Value = 744.7778
$ws.Range("K122").Value = 6703.000199999999
$ws.Range("M122").Value = -4253.000199999999
$ws.Range("H131").Value = 835.71
$ws.Range("I131").Value = 356.15384
$ws.Range("J131").Value = 907.3678
$ws.Range("K131").Value = 1068.46152
$ws.Range("L131").Value = 2722.1034
$ws.Range("M131").Value = 3971.53848
$ws.Range("N131").Value = -12802.1034
$ws.Range("H135").Value = 651.53845
$ws.Range("I135").Value = 541.44446
$ws.Range("J135").Value = 899.25
$ws.Range("K135").Value = 4873.00014
$ws.Range("L135").Value = 8093.25
$ws.Range("M135").Value = -2338.00014
$ws.Range("N135").Value = -13163.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 381.5
$ws.Range("J107").Value = 494.5
$ws.Range("L107").Value = 494.5
$ws.Range("N107").Value = -4334.5
$ws.Range("H123").Value = 17241.666
$ws.Range("J123").Value = 17241.666
$ws.Range("L123").Value = 17241.666
$ws.Range("N123").Value = -22141.666
$ws.Range("H126").Value = 4715706
$ws.Range("I126").Value = 18522966
$ws.Range("J126").Value = 113285.89
$ws.Range("K126").Value = 55568898
$ws.Range("L126").Value = 339857.67
$ws.Range("M126").Value = -55566428
$ws.Range("N126").Value = -344797.67
$ws.Range("H132").Value = 1101511.5
$ws.Range("I132").Value = 1750039
$ws.Range("K132").Value = 5250117
$ws.Range("M132").Value = -5247587

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 402000
$ws.Range("J2").Value = 10000
$ws.Range("L2").Value = 10000
$ws.Range("N2").Value = -10224
$ws.Range("H122").Value = 4154.6665
$ws.Range("I122").Value = 4686.1113
$ws.Range("J122").Value = 3623.2222
$ws.Range("K122").Value = 14058.3339
$ws.Range("L122").Value = 10869.6666
$ws.Range("M122").Value = -11608.3339
$ws.Range("N122").Value = -15769.6666
$ws.Range("H132").Value = 1764.0513
$ws.Range("I132").Value = 1132.7878
$ws.Range("J132").Value = 5236
$ws.Range("K132").Value = 3398.3634
$ws.Range("L132").Value = 15708
$ws.Range("M132").Value = -868.3634000000002
$ws.Range("N132").Value = -20768
$ws.Range("H136").Value = 2000.5869
$ws.Range("I136").Value = 1229.7142
$ws.Range("K136").Value = 3689.1426
$ws.Range("M136").Value = -1139.1426

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1711.9788
$ws.Range("I132").Value = 1122.3513
$ws.Range("K132").Value = 3367.0539
$ws.Range("M132").Value = -837.0538999999999
$ws.Range("H135").Value = 85059
$ws.Range("J135").Value = 85059
$ws.Range("L135").Value = 85059
$ws.Range("N135").Value = -95199
$ws.Range("H136").Value = 23149894
$ws.Range("I136").Value = 29241092
$ws.Range("J136").Value = 3339.8
$ws.Range("K136").Value = 87723276
$ws.Range("L136").Value = 10019.4
$ws.Range("M136").Value = -87720726
$ws.Range("N136").Value = -15119.4
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = $null
